$wb = $excel.ActiveWorkbook

function Remove-HyperlinkAt($ws, $addr) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address(0, 0) -eq $addr) {
            $h.Delete()
        }
    }
}

# --- Overview sheet: status text changes from "Ready for handoff" to "Handoff transform failed"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

# --- zh-cn sheet: handoff attempt reset (transform failed, nothing handed off)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = "Handoff transform failed"

Remove-HyperlinkAt $wsZh "C2"
$wsZh.Range("C2").Clear()

$wsZh.Range("D2").Value = "0001-01-01 00:00:00"
$wsZh.Range("G2").Value = "0001-01-01 00:00:00"
$wsZh.Range("H2").Value = "Ignored"

$wsZh.Range("D3").Value = "0001-01-01 00:00:00"
$wsZh.Range("G3").Value = "0001-01-01 00:00:00"
$wsZh.Range("H3").Value = "Ignored"

# --- de-de sheet: same reset
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = "Handoff transform failed"

Remove-HyperlinkAt $wsDe "C2"
$wsDe.Range("C2").Clear()

$wsDe.Range("D2").Value = "0001-01-01 00:00:00"
$wsDe.Range("G2").Value = "0001-01-01 00:00:00"
$wsDe.Range("H2").Value = "Ignored"

$wsDe.Range("D3").Value = "0001-01-01 00:00:00"
$wsDe.Range("G3").Value = "0001-01-01 00:00:00"
$wsDe.Range("H3").Value = "Ignored"
